# #5: cash & deposit done
# Rebuild the "存款" (bank deposit) sheet with bank/deposit_type/currency headers
# plus the full legislator metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) that the other property
# sheets (股票, 基金受益憑證, ...) already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

$header = @("bank", "deposit_type", "currency", "owner", "total", "property_category", "category", "date", "legislator_name", "legislator_id", "source_file", "index")
for ($i = 0; $i -lt $header.Count; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $header[$i]
}

$rows = @(
    @(46, "臺灣銀行群賢分行",       "活期存款", "新臺幣", "田秋堇", 97097,      "deposit", "normal", "2012-04-10", "田秋堇", 1316, "tmp9b251", 46),
    @(47, "台北富邦商業銀行羅東分行", "活期存款", "新臺幣", "田秋堇", 66310,      "deposit", "normal", "2012-04-10", "田秋堇", 1316, "tmp9b251", 47),
    @(48, "台北富邦商業銀行羅東分行", "活期存款", "澳幣",   "田秋堇", 4463.54,    "deposit", "normal", "2012-04-10", "田秋堇", 1316, "tmp9b251", 48),
    @(49, "台北富邦商業銀行羅東分行", "活期存款", "美金",   "田秋堇", 1073920.62, "deposit", "normal", "2012-04-10", "田秋堇", 1316, "tmp9b251", 49),
    @(50, "台北富邦商業銀行羅東分行", "活期存款", "新臺幣", "劉守成", 169017,     "deposit", "normal", "2012-04-10", "田秋堇", 1316, "tmp9b251", 50),
    @(51, "台北富邦商業銀行羅東分行", "活期存款", "曰圓",   "劉守成", 3227.88,    "deposit", "normal", "2012-04-10", "田秋堇", 1316, "tmp9b251", 51),
    @(52, "台北富邦商業銀行羅東分行", "活期存款", "歐元",   "劉守成", 0.77,       "deposit", "normal", "2012-04-10", "田秋堇", 1316, "tmp9b251", 52)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $row = $rows[$r]
    $excelRow = 2 + $r
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($excelRow, 1 + $c).Value = $row[$c]
    }
}
